# Cambios en asignacion inicial, recotizada, duplicados.
# Cambio en clase load para incluir posibles tables o schemas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Fill in missing "status_name" values (asignacion inicial / recotizada / duplicados) ---
$a13 = $ws.Range("A13").Value()
$a23 = $ws.Range("A23").Value()
$a24 = $ws.Range("A24").Value()
$ws.Range("B13").Value = $a13
$ws.Range("B23").Value = $a23
$ws.Range("B24").Value = $a24

# --- Add a new "plataforma" column to the table, filled with "vicidial" ---
$tbl.ListColumns.Add() | Out-Null

# Carry over the formatting from the last existing column so the new
# header/data cells match the rest of the table.
$ws.Range("J1:J24").Copy() | Out-Null
$ws.Range("K1:K24").PasteSpecial(-4122) | Out-Null

$ws.Range("K1").Value = "plataforma"
$ws.Range("K2:K24").Value = "vicidial"

$ws.Range("B13").Select() | Out-Null
